$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.161.78'
$ws.Range('D3').Value = '1.563.03'
$ws.Range('E3').Value = '  -1.60%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '206.71'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('E6').Value = '  -1.50%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.248'
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0592'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.12%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0861'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.96%  '
$ws.Range('D12').Value = '1.784.26'
$ws.Range('E12').Value = '  -1.65%  '
$ws.Range('D13').Value = '1.569.04'
$ws.Range('E13').Value = '  -1.12%  '
$ws.Range('E14').Value = '  -2.13%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.516'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.61%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '62.93'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.79%  '
$ws.Range('D17').Value = '27.155.27'
$ws.Range('E17').Value = '  -1.84%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '213.48'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.69%  '
$ws.Range('E19').Value = '  -1.11%  '
$ws.Range('E20').Value = '  -1.37%  '
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('E22').Value = '  -0.52%  '
$ws.Range('E23').Value = '  -2.15%  '
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '152.17'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.84%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.57'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -3.83%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '14.88'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.62%  '
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('E29').Value = '  -1.49%  '
$ws.Range('E30').Value = '  -1.25%  '
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.17'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.77%  '
$ws.Range('D33').Value = '1.383.02'
$ws.Range('E33').Value = '  +0.93%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.95'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.63%  '
$ws.Range('E35').Value = '  +0.57%  '
$ws.Range('E36').Value = '  -1.08%  '
$ws.Range('E37').Value = '  -4.15%  '
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.815'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.19%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.517'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.54%  '
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.992'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.98%  '
$ws.Range('E43').Value = '  +2.88%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '63.46'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.99%  '
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').Value = '1.697.05'
$ws.Range('E47').Value = '  -1.65%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '85.60'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.22%  '
$ws.Range('D49').Value = '0.0₇0998'
$ws.Range('E49').Value = '  -0.86%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0492'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.67%  '
$ws.Range('E51').Value = '  -0.07%  '
